$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix comma -> period typos in two "Razon social"/"Nombre Fantasia" entries
$ws.Range("E25").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F25").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E68").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F68").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E48").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# Normalize "Importe" column from Argentine decimal notation (1.234,56)
# to plain-dot decimal notation (1234.56), keeping the values as text
# so the trailing zeros / exact digits survive (match as typed by scraper).
$importeRange = $ws.Range("H2:H114")
$importeRange.NumberFormat = "@"

$ws.Range("H2").Value = "540.00"
$ws.Range("H3").Value = "24999.00"
$ws.Range("H4").Value = "514.00"
$ws.Range("H5").Value = "440.00"
$ws.Range("H6").Value = "90.00"
$ws.Range("H7").Value = "208.00"
$ws.Range("H8").Value = "26110.00"
$ws.Range("H9").Value = "66894.76"
$ws.Range("H10").Value = "9147.55"
$ws.Range("H11").Value = "7482.92"
$ws.Range("H12").Value = "4083.95"
$ws.Range("H13").Value = "6183.09"
$ws.Range("H14").Value = "270.60"
$ws.Range("H15").Value = "6724.30"
$ws.Range("H16").Value = "4630.00"
$ws.Range("H17").Value = "30.00"
$ws.Range("H18").Value = "1730.00"
$ws.Range("H19").Value = "2934.34"
$ws.Range("H20").Value = "19455.00"
$ws.Range("H21").Value = "1786.00"
$ws.Range("H22").Value = "962.00"
$ws.Range("H23").Value = "24345.00"
$ws.Range("H24").Value = "1850.00"
$ws.Range("H25").Value = "198.17"
$ws.Range("H26").Value = "1803.54"
$ws.Range("H27").Value = "266.80"
$ws.Range("H28").Value = "1040.27"
$ws.Range("H29").Value = "248.00"
$ws.Range("H30").Value = "17636.00"
$ws.Range("H31").Value = "452.46"
$ws.Range("H32").Value = "279.00"
$ws.Range("H33").Value = "522.23"
$ws.Range("H34").Value = "251.08"
$ws.Range("H35").Value = "542.00"
$ws.Range("H36").Value = "270.00"
$ws.Range("H37").Value = "15680.00"
$ws.Range("H38").Value = "6368.00"
$ws.Range("H39").Value = "30.00"
$ws.Range("H40").Value = "5440.00"
$ws.Range("H41").Value = "170.00"
$ws.Range("H42").Value = "1700.00"
$ws.Range("H43").Value = "1305.00"
$ws.Range("H44").Value = "316.60"
$ws.Range("H45").Value = "12194.80"
$ws.Range("H46").Value = "294.00"
$ws.Range("H47").Value = "1164.12"
$ws.Range("H48").Value = "390.00"
$ws.Range("H49").Value = "345.00"
$ws.Range("H50").Value = "180.00"
$ws.Range("H51").Value = "4734.00"
$ws.Range("H52").Value = "0.11"
$ws.Range("H53").Value = "0.35"
$ws.Range("H54").Value = "100171.00"
$ws.Range("H55").Value = "0.70"
$ws.Range("H56").Value = "38.80"
$ws.Range("H57").Value = "433.00"
$ws.Range("H58").Value = "2104.20"
$ws.Range("H59").Value = "14.95"
$ws.Range("H60").Value = "31.04"
$ws.Range("H61").Value = "90.00"
$ws.Range("H62").Value = "3090.00"
$ws.Range("H63").Value = "191.00"
$ws.Range("H64").Value = "230.00"
$ws.Range("H65").Value = "336.00"
$ws.Range("H66").Value = "7179.00"
$ws.Range("H67").Value = "36.00"
$ws.Range("H68").Value = "737.76"
$ws.Range("H69").Value = "33.00"
$ws.Range("H70").Value = "400.00"
$ws.Range("H71").Value = "785.00"
$ws.Range("H72").Value = "782.86"
$ws.Range("H73").Value = "144.93"
$ws.Range("H74").Value = "99.00"
$ws.Range("H75").Value = "199.00"
$ws.Range("H76").Value = "39.65"
$ws.Range("H77").Value = "1053.25"
$ws.Range("H78").Value = "5000.00"
$ws.Range("H79").Value = "800.00"
$ws.Range("H80").Value = "2846.94"
$ws.Range("H81").Value = "1200.00"
$ws.Range("H82").Value = "13504.58"
$ws.Range("H83").Value = "297.00"
$ws.Range("H84").Value = "90.00"
$ws.Range("H85").Value = "1850.00"
$ws.Range("H86").Value = "350.00"
$ws.Range("H87").Value = "84.00"
$ws.Range("H88").Value = "195.00"
$ws.Range("H89").Value = "251.62"
$ws.Range("H90").Value = "426.00"
$ws.Range("H91").Value = "817.62"
$ws.Range("H92").Value = "29897.02"
$ws.Range("H93").Value = "240.00"
$ws.Range("H94").Value = "2.00"
$ws.Range("H95").Value = "488.00"
$ws.Range("H96").Value = "730.00"
$ws.Range("H97").Value = "2039.00"
$ws.Range("H98").Value = "450.00"
$ws.Range("H99").Value = "1150.00"
$ws.Range("H100").Value = "1000.00"
$ws.Range("H101").Value = "660.00"
$ws.Range("H102").Value = "14981.00"
$ws.Range("H103").Value = "138.00"
$ws.Range("H104").Value = "480.00"
$ws.Range("H105").Value = "37538.95"
$ws.Range("H106").Value = "2704.13"
$ws.Range("H107").Value = "40507.78"
$ws.Range("H108").Value = "850.00"
$ws.Range("H109").Value = "8469.01"
$ws.Range("H110").Value = "309857.23"
$ws.Range("H111").Value = "44000.00"
$ws.Range("H112").Value = "2559.28"
$ws.Range("H113").Value = "380.00"
$ws.Range("H114").Value = "350.00"
